$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.580.30"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "2.029.51"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "226.44"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").Value = "0.606"
$ws.Range("E6").Value = "  -1.02%  "

$ws.Range("D7").Value = "59.62"
$ws.Range("E7").Value = "  -1.05%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("E10").Value = "  +2.08%  "

$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").Value = "2.333.18"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "14.38"
$ws.Range("E13").Value = "  -1.57%  "

$ws.Range("D14").Value = "21.05"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").Value = "5.46"
$ws.Range("E15").Value = "  +4.30%  "

$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("D17").Value = "2.033.52"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "37.546.89"
$ws.Range("E18").Value = "  -0.73%  "

$ws.Range("E19").Value = "  -2.29%  "

$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").Value = "223.56"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("E24").Value = "  -0.22%  "

$ws.Range("E25").Value = "  +3.13%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.39"
$ws.Range("E26").Value = "  +2.50%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "167.44"
$ws.Range("E27").Value = "  +1.43%  "

$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").Value = "18.71"
$ws.Range("E29").Value = "  -0.91%  "

$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("E32").Value = "  +8.97%  "

$ws.Range("E33").Value = "  -1.85%  "

$ws.Range("D34").Value = "0.0603"
$ws.Range("E34").Value = "  +0.09%  "

$ws.Range("D35").Value = "4.45"
$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("D36").Value = "6.48"
$ws.Range("E36").Value = "  +2.76%  "

$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  +3.51%  "

$ws.Range("E38").Value = "  +5.15%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").Value = "17.92"
$ws.Range("E40").Value = "  +8.29%  "

$ws.Range("D41").Value = "1.522.88"

$ws.Range("D42").Value = "96.82"
$ws.Range("E42").Value = "  -0.26%  "

$ws.Range("E43").Value = "  -1.12%  "

$ws.Range("E44").Value = "  +1.31%  "

$ws.Range("D45").Value = "0.0906"
$ws.Range("E45").Value = "  -1.27%  "

$ws.Range("D46").Value = "4.18"
$ws.Range("E46").Value = "  +4.55%  "

$ws.Range("E47").Value = "  -0.61%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("E49").Value = "  -1.05%  "

$ws.Range("D50").Value = "7.04"
$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").Value = "2.223.29"
$ws.Range("E51").Value = "  +0.03%  "
